$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.456404
$ws.Range("H2").Value = 1.369212
$ws.Range("I2").Value = 0.01914960767004715
$ws.Range("J2").Value = 0.01914960767004715
$ws.Range("M2").Value = 1.116695
$ws.Range("N2").Value = 3.350085
$ws.Range("O2").Value = 0.008174214292497491
$ws.Range("P2").Value = 0.008174214292497492
$ws.Range("Q2").Value = 0.50966406478
$ws.Range("R2").Value = 4.58697658302
$ws.Range("S2").Value = 0.000156532996712219
$ws.Range("T2").Value = 0.000156532996712219
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.456404
$ws.Range("H3").Value = 1.369212
$ws.Range("I3").Value = 0.01914960767004715
$ws.Range("J3").Value = 0.01914960767004715
$ws.Range("O3").Value = 0.8193429796700005
$ws.Range("P3").Value = 0.8193429796700005
$ws.Range("Q3").Value = 51.08621557069333
$ws.Range("R3").Value = 459.77594013624
$ws.Range("S3").Value = 0.01569009660788793
$ws.Range("T3").Value = 0.01569009660788793
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.456404
$ws.Range("H4").Value = 1.369212
$ws.Range("I4").Value = 0.01914960767004715
$ws.Range("J4").Value = 0.01914960767004715
$ws.Range("O4").Value = 0.172482806037502
$ws.Range("P4").Value = 0.1724828060375021
$ws.Range("Q4").Value = 10.75434101482
$ws.Range("R4").Value = 96.78906913338001
$ws.Range("S4").Value = 0.003302978065447004
$ws.Range("T4").Value = 0.003302978065447005
$ws.Range("I5").Value = 0.8285024587002443
$ws.Range("J5").Value = 0.8285024587002443
$ws.Range("M5").Value = 1.116695
$ws.Range("N5").Value = 3.350085
$ws.Range("O5").Value = 0.008174214292497491
$ws.Range("P5").Value = 0.008174214292497492
$ws.Range("Q5").Value = 22.05047424767167
$ws.Range("R5").Value = 198.454268229045
$ws.Range("S5").Value = 0.006772356639276849
$ws.Range("T5").Value = 0.00677235663927685
$ws.Range("I6").Value = 0.8285024587002443
$ws.Range("J6").Value = 0.8285024587002443
$ws.Range("O6").Value = 0.8193429796700005
$ws.Range("P6").Value = 0.8193429796700005
$ws.Range("S6").Value = 0.6788276731753796
$ws.Range("T6").Value = 0.6788276731753796
$ws.Range("I7").Value = 0.8285024587002443
$ws.Range("J7").Value = 0.8285024587002443
$ws.Range("O7").Value = 0.172482806037502
$ws.Range("P7").Value = 0.1724828060375021
$ws.Range("S7").Value = 0.1429024288855878
$ws.Range("T7").Value = 0.1429024288855878
$ws.Range("I8").Value = 0.1523479336297086
$ws.Range("J8").Value = 0.1523479336297086
$ws.Range("M8").Value = 1.116695
$ws.Range("N8").Value = 3.350085
$ws.Range("O8").Value = 0.008174214292497491
$ws.Range("P8").Value = 0.008174214292497492
$ws.Range("Q8").Value = 4.054718428305001
$ws.Range("R8").Value = 36.49246585474501
$ws.Range("S8").Value = 0.001245324656508423
$ws.Range("T8").Value = 0.001245324656508424
$ws.Range("I9").Value = 0.1523479336297086
$ws.Range("J9").Value = 0.1523479336297086
$ws.Range("O9").Value = 0.8193429796700005
$ws.Range("P9").Value = 0.8193429796700005
$ws.Range("S9").Value = 0.1248252098867329
$ws.Range("T9").Value = 0.1248252098867329
$ws.Range("I10").Value = 0.1523479336297086
$ws.Range("J10").Value = 0.1523479336297086
$ws.Range("O10").Value = 0.172482806037502
$ws.Range("P10").Value = 0.1724828060375021
$ws.Range("S10").Value = 0.02627739908646727
$ws.Range("T10").Value = 0.02627739908646727
